{"js": "// Fix the duplicated \"s\u00fahvezdie\" word: the paragraph currently reads\n// \"...pozorova\u0165 s\u00fahvezdie S\u00fahvezdie Herkules...\" and should read\n// \"...pozorova\u0165 S\u00fahvezdie Herkules...\" (remove the stray lowercase\n// duplicate, keep the capitalized \"S\u00fahvezdie Herkules\").\nconst results = context.document.body.search(\"s\u00fahvezdie S\u00fahvezdie Herkules\", {\n  matchCase: true\n});\nresults.load(\"items\");\nawait context.sync();\n\nfor (let i = 0; i < results.items.length; i++) {\n  results.items[i].insertText(\"S\u00fahvezdie Herkules\", Word.InsertLocation.replace);\n}\nawait context.sync();\n", "ps1": "# Fix the duplicated \"s\u00fahvezdie\" word: the paragraph currently reads\n# \"...pozorova\u0165 s\u00fahvezdie S\u00fahvezdie Herkules...\" and should read\n# \"...pozorova\u0165 S\u00fahvezdie Herkules...\" (remove the stray lowercase\n# duplicate, keep the capitalized \"S\u00fahvezdie Herkules\"). There are four\n# identical occurrences in the document, so replace all of them.\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"s\u00fahvezdie S\u00fahvezdie Herkules\"\n$find.Replacement.Text = \"S\u00fahvezdie Herkules\"\n$find.Forward = $true\n$find.Wrap = 1          # wdFindContinue\n$find.Format = $false\n$find.MatchCase = $true\n$find.MatchWholeWord = $false\n$find.MatchWildcards = $false\n\n$find.Execute([ref]$find.Text, [ref]$find.MatchCase, [ref]$find.MatchWholeWord, `\n  [ref]$find.MatchWildcards, $null, $null, [ref]$find.Forward, [ref]$find.Wrap, `\n  $null, [ref]$find.Replacement.Text, 2)   # 2 = wdReplaceAll\n"}
